$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Word Count")

# Add new Battery SOC table entry: date 2024-06-23 (serial 45466) with 885 words
$ws.Range("A11").Value = 45466
$ws.Range("B11").Value = 885

# Update the selection to match the saved state (cell A12 selected)
$ws.Range("A12").Select()
